$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 to the new (deeper) URL, preserving D4's existing value.
$ws.Range("B3").Value = "https://development1.advantageclub.co/in/rewards/home"

# Add new rows with login credentials used by the automation script.
$ws.Range("B11").Value = "dheerajc@advantageclub.in"
$ws.Range("B12").Value = "Dheeraj@4321"
$ws.Range("B15").Value = "adminadvantage"
$ws.Range("B16").Value = "@Advantage_1"

# Turn B11, B12 and B3 (in that order) into hyperlinks - mirrors rId1/rId2/rId3 order.
$ws.Hyperlinks.Add($ws.Range("B11"), "dheerajc@advantageclub.in")
$ws.Hyperlinks.Add($ws.Range("B12"), "Dheeraj@4321")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://development1.advantageclub.co/in/rewards/home")

# Widen column B to fit the longer strings, matching the new layout
# (best-fit width for this content renders as 55 characters wide).
$ws.Columns("B").ColumnWidth = 55 - 5/6

# Update the visible selection to B6, as captured in the saved workbook.
$ws.Range("B6").Select()
